$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H88").Value = 6666.3335
$ws_ALC.Range("I88").Value = 4999.5
$ws_ALC.Range("J88").Value = 10000
$ws_ALC.Range("K88").Value = 4999.5
$ws_ALC.Range("L88").Value = 10000
$ws_ALC.Range("M88").Value = -4593.5
$ws_ALC.Range("N88").Value = -10812
$ws_ALC.Range("H91").Value = 6666.3335
$ws_ALC.Range("I91").Value = 4999.5
$ws_ALC.Range("J91").Value = 10000
$ws_ALC.Range("K91").Value = 4999.5
$ws_ALC.Range("L91").Value = 10000
$ws_ALC.Range("M91").Value = -3595.5
$ws_ALC.Range("N91").Value = -12808
$ws_ALC.Range("H101").Value = 1433.3334
$ws_ALC.Range("J101").Value = 2000
$ws_ALC.Range("L101").Value = 6000
$ws_ALC.Range("N101").Value = -9244
$ws_ALC.Range("H132").Value = 1299.4839
$ws_ALC.Range("I132").Value = 1209.4667
$ws_ALC.Range("K132").Value = 3628.4001
$ws_ALC.Range("M132").Value = -1098.4001
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 2991.3438
$ws_ARM.Range("I32").Value = 1827.94
$ws_ARM.Range("K32").Value = 1827.94
$ws_ARM.Range("M32").Value = -1540.94
$ws_ARM.Range("H61").Value = 4083.7778
$ws_ARM.Range("J61").Value = 5375.1113
$ws_ARM.Range("L61").Value = 5375.1113
$ws_ARM.Range("N61").Value = -5799.1113
$ws_ARM.Range("H74").Value = 4728.5625
$ws_ARM.Range("I74").Value = 4714.643
$ws_ARM.Range("K74").Value = 4714.643
$ws_ARM.Range("M74").Value = -3840.643
$ws_ARM.Range("H77").Value = 4728.5625
$ws_ARM.Range("I77").Value = 4714.643
$ws_ARM.Range("K77").Value = 23573.215
$ws_ARM.Range("M77").Value = -19205.215
$ws_ARM.Range("H132").Value = 2717.7144
$ws_ARM.Range("I132").Value = 1964.6666
$ws_ARM.Range("J132").Value = 3721.7778
$ws_ARM.Range("K132").Value = 5893.9998
$ws_ARM.Range("L132").Value = 11165.3334
$ws_ARM.Range("M132").Value = -3363.9998
$ws_ARM.Range("N132").Value = -16225.3334
$ws_ARM.Range("H136").Value = 4083.7778
$ws_ARM.Range("J136").Value = 5375.1113
$ws_ARM.Range("L136").Value = 16125.3339
$ws_ARM.Range("N136").Value = -21225.3339
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H7").Value = 3000
$ws_BSM.Range("I7").Value = 3000
$ws_BSM.Range("K7").Value = 3000
$ws_BSM.Range("M7").Value = -2887
$ws_BSM.Range("H86").Value = 202499.6
$ws_BSM.Range("I86").Value = 2624.5
$ws_BSM.Range("K86").Value = 2624.5
$ws_BSM.Range("M86").Value = -1501.5
$ws_BSM.Range("H89").Value = 202499.6
$ws_BSM.Range("I89").Value = 2624.5
$ws_BSM.Range("K89").Value = 13122.5
$ws_BSM.Range("M89").Value = -7506.5
$ws_BSM.Range("H134").Value = 7603.9033
$ws_BSM.Range("I134").Value = 7774.36
$ws_BSM.Range("K134").Value = 23323.08
$ws_BSM.Range("M134").Value = -20788.08
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H4").Value = 9500
$ws_CRP.Range("J4").Value = 9500
$ws_CRP.Range("L4").Value = 9500
$ws_CRP.Range("N4").Value = -9724
$ws_CRP.Range("H31").Value = 2467.795
$ws_CRP.Range("I31").Value = 1060.56
$ws_CRP.Range("K31").Value = 1060.56
$ws_CRP.Range("M31").Value = -765.5599999999999
$ws_CRP.Range("H34").Value = 2467.795
$ws_CRP.Range("I34").Value = 1060.56
$ws_CRP.Range("K34").Value = 1060.56
$ws_CRP.Range("M34").Value = -858.5599999999999
$ws_CRP.Range("H99").Value = 2285.8
$ws_CRP.Range("I99").Value = 1739.5
$ws_CRP.Range("J99").Value = 2650
$ws_CRP.Range("K99").Value = 1739.5
$ws_CRP.Range("L99").Value = 2650
$ws_CRP.Range("M99").Value = -241.5
$ws_CRP.Range("N99").Value = -5646
$ws_CRP.Range("H126").Value = 2285.8
$ws_CRP.Range("I126").Value = 1739.5
$ws_CRP.Range("J126").Value = 2650
$ws_CRP.Range("K126").Value = 5218.5
$ws_CRP.Range("L126").Value = 7950
$ws_CRP.Range("M126").Value = -2748.5
$ws_CRP.Range("N126").Value = -12890
$ws_CRP.Range("H134").Value = 1246.5454
$ws_CRP.Range("I134").Value = 1249.875
$ws_CRP.Range("K134").Value = 3749.625
$ws_CRP.Range("M134").Value = -1214.625
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 20849.139
$ws_CUL.Range("J4").Value = 373484.5
$ws_CUL.Range("L4").Value = 1120453.5
$ws_CUL.Range("N4").Value = -1120677.5
$ws_CUL.Range("H80").Value = 3125
$ws_CUL.Range("I80").Value = 5000
$ws_CUL.Range("J80").Value = 2500
$ws_CUL.Range("K80").Value = 15000
$ws_CUL.Range("L80").Value = 7500
$ws_CUL.Range("M80").Value = -14064
$ws_CUL.Range("N80").Value = -9372
$ws_CUL.Range("H83").Value = 3125
$ws_CUL.Range("I83").Value = 5000
$ws_CUL.Range("J83").Value = 2500
$ws_CUL.Range("K83").Value = 45000
$ws_CUL.Range("L83").Value = 22500
$ws_CUL.Range("M83").Value = -40320
$ws_CUL.Range("N83").Value = -31860
$ws_CUL.Range("H87").Value = 11924.714
$ws_CUL.Range("I87").Value = 1157.6666
$ws_CUL.Range("J87").Value = 20000
$ws_CUL.Range("K87").Value = 3472.9998
$ws_CUL.Range("L87").Value = 60000
$ws_CUL.Range("M87").Value = -2224.9998
$ws_CUL.Range("N87").Value = -62496
$ws_CUL.Range("H90").Value = 11924.714
$ws_CUL.Range("I90").Value = 1157.6666
$ws_CUL.Range("J90").Value = 20000
$ws_CUL.Range("K90").Value = 10418.9994
$ws_CUL.Range("L90").Value = 180000
$ws_CUL.Range("M90").Value = -4178.999400000001
$ws_CUL.Range("N90").Value = -192480
$ws_CUL.Range("H103").Value = 1284.5
$ws_CUL.Range("I103").Value = 929.8333
$ws_CUL.Range("K103").Value = 2789.4999
$ws_CUL.Range("M103").Value = -1910.4999
$ws_CUL.Range("H129").Value = 32381.39
$ws_CUL.Range("J129").Value = 46215.375
$ws_CUL.Range("L129").Value = 138646.125
$ws_CUL.Range("N129").Value = -148646.125
$ws_CUL.Range("H131").Value = 7949907
$ws_CUL.Range("J131").Value = 14756.947
$ws_CUL.Range("L131").Value = 44270.841
$ws_CUL.Range("N131").Value = -54350.841
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 1771.84
$ws_GSM.Range("I102").Value = 1660.3334
$ws_GSM.Range("J102").Value = 1874.7693
$ws_GSM.Range("K102").Value = 1660.3334
$ws_GSM.Range("L102").Value = 1874.7693
$ws_GSM.Range("M102").Value = -38.33339999999998
$ws_GSM.Range("N102").Value = -5118.7693
$ws_GSM.Range("H122").Value = 1769.6428
$ws_GSM.Range("I122").Value = 1684.7059
$ws_GSM.Range("J122").Value = 1900.909
$ws_GSM.Range("K122").Value = 5054.1177
$ws_GSM.Range("L122").Value = 5702.727000000001
$ws_GSM.Range("M122").Value = -2604.1177
$ws_GSM.Range("N122").Value = -10602.727
$ws_GSM.Range("H132").Value = 2977.1924
$ws_GSM.Range("I132").Value = 2204.889
$ws_GSM.Range("J132").Value = 4714.875
$ws_GSM.Range("K132").Value = 6614.667
$ws_GSM.Range("L132").Value = 14144.625
$ws_GSM.Range("M132").Value = -4084.667
$ws_GSM.Range("N132").Value = -19204.625
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 4916.2856
$ws_LTW.Range("I7").Value = 2594.875
$ws_LTW.Range("J7").Value = 6344.846
$ws_LTW.Range("K7").Value = 2594.875
$ws_LTW.Range("L7").Value = 6344.846
$ws_LTW.Range("M7").Value = -2482.875
$ws_LTW.Range("N7").Value = -6568.846
$ws_LTW.Range("H126").Value = 4916.2856
$ws_LTW.Range("I126").Value = 2594.875
$ws_LTW.Range("J126").Value = 6344.846
$ws_LTW.Range("K126").Value = 7784.625
$ws_LTW.Range("L126").Value = 19034.538
$ws_LTW.Range("M126").Value = -5314.625
$ws_LTW.Range("N126").Value = -23974.538
$ws_LTW.Range("H132").Value = 1947
$ws_LTW.Range("I132").Value = 1450.2858
$ws_LTW.Range("K132").Value = 4350.857400000001
$ws_LTW.Range("M132").Value = -1820.857400000001
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H25").Value = 15000
$ws_WVR.Range("J25").Value = 15000
$ws_WVR.Range("L25").Value = 15000
$ws_WVR.Range("N25").Value = -15586
$ws_WVR.Range("H132").Value = 5899.636
$ws_WVR.Range("I132").Value = 1199.5
$ws_WVR.Range("K132").Value = 3598.5
$ws_WVR.Range("M132").Value = -1068.5
